# "first run of Vlookup"
# Populate the "merged_data_python" summary sheet's "Points Per Game" block
# (columns B:F, rows 2:33) by looking each team up in the "Points Per Game"
# sheet, then flatten the formulas down to their computed values (matching
# a paste-values workflow after a first VLOOKUP pass).

$wb = $excel.ActiveWorkbook

$srcName = "Points Per Game"
$dstName = "merged_data_python"

$src = $wb.Worksheets.Item($srcName)
$dst = $wb.Worksheets.Item($dstName)

# VLOOKUP each row's team name (column A) against the source sheet, pulling
# back the matching 2023 / Last 3 / Last 1 / Home / Away columns (B:F).
$dst.Range("B2:F33").Formula = "=VLOOKUP(`$A2,'Points Per Game'!`$A:`$F,COLUMN(),FALSE)"

# Freeze the lookups into plain values, same as the committed workbook
# (no live formulas remain in the saved file).
$dst.Range("B2:F33").Copy()
$dst.Range("B2:F33").PasteSpecial(-4163)
